# Append the 2025-10-02 profit allocation row after running the model.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 31

# Column A holds the date as literal text (like the rows above it), not an
# Excel date serial, so force text format before assigning the value and
# then drop back to the sheet's normal style to avoid leaving a stray
# per-cell format behind.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "10/02/2025"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 0.1380907263517728
$ws.Cells.Item($row, 3).Value = 0.8619092736482272
